# Updated symbol list on Mon Jan  9 13:24:41 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume/hour snapshot in Sheet1 (rows 2-51):
#   - Column D ("Price")      updated for rows with live market data
#   - Column E ("Volume(1h)") updated for rows with live market data
#   - Column G ("Hora")       bumped from 12 -> 13 for every data row
#
# D/E/G are stored as literal text in the workbook (e.g. "12", "6.66%"),
# not as numbers/percentages, so we force the Text number format before
# writing the values (otherwise Excel would silently reinterpret "13" as
# a number or "6.79%" as a 0.0679 percentage). The format is reset back
# to the default "Normal" style afterwards so no stray formatting is left
# behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("D2:G51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "279.29"
$ws.Range("E2").Value = "6.79%"
$ws.Range("G2").Value = "13"
$ws.Range("D3").Value = "27.39"
$ws.Range("E3").Value = "2.69%"
$ws.Range("G3").Value = "13"
$ws.Range("D4").Value = "4.812"
$ws.Range("E4").Value = "2.37%"
$ws.Range("G4").Value = "13"
$ws.Range("D5").Value = "0.06344"
$ws.Range("E5").Value = "4.35%"
$ws.Range("G5").Value = "13"
$ws.Range("D6").Value = "6.946"
$ws.Range("E6").Value = "3.66%"
$ws.Range("G6").Value = "13"
$ws.Range("D7").Value = "3.372"
$ws.Range("E7").Value = "6.40%"
$ws.Range("G7").Value = "13"
$ws.Range("D8").Value = "0.8810"
$ws.Range("E8").Value = "3.50%"
$ws.Range("G8").Value = "13"
$ws.Range("D9").Value = "0.9576"
$ws.Range("E9").Value = "4.99%"
$ws.Range("G9").Value = "13"
$ws.Range("D10").Value = "0.1483"
$ws.Range("E10").Value = "5.66%"
$ws.Range("G10").Value = "13"
$ws.Range("D11").Value = "0.05187"
$ws.Range("E11").Value = "0.90%"
$ws.Range("G11").Value = "13"
$ws.Range("E12").Value = "2.65%"
$ws.Range("G12").Value = "13"
$ws.Range("D13").Value = "0.03141"
$ws.Range("E13").Value = "0.88%"
$ws.Range("G13").Value = "13"
$ws.Range("D14").Value = "0.09065"
$ws.Range("E14").Value = "0.33%"
$ws.Range("G14").Value = "13"
$ws.Range("D15").Value = "0.001567"
$ws.Range("E15").Value = "2.04%"
$ws.Range("G15").Value = "13"
$ws.Range("D16").Value = "0.0006259"
$ws.Range("E16").Value = "1.16%"
$ws.Range("G16").Value = "13"
$ws.Range("D17").Value = "0.005875"
$ws.Range("E17").Value = "-3.49%"
$ws.Range("G17").Value = "13"
$ws.Range("D18").Value = "3.463"
$ws.Range("E18").Value = "0.42%"
$ws.Range("G18").Value = "13"
$ws.Range("D19").Value = "2.297"
$ws.Range("E19").Value = "6.01%"
$ws.Range("G19").Value = "13"
$ws.Range("D20").Value = "0.3122"
$ws.Range("E20").Value = "1.61%"
$ws.Range("G20").Value = "13"
$ws.Range("D21").Value = "0.1291"
$ws.Range("E21").Value = "-0.72%"
$ws.Range("G21").Value = "13"
$ws.Range("D22").Value = "3.865"
$ws.Range("E22").Value = "-5.97%"
$ws.Range("G22").Value = "13"
$ws.Range("D23").Value = "0.04318"
$ws.Range("E23").Value = "2.15%"
$ws.Range("G23").Value = "13"
$ws.Range("E24").Value = "-0.28%"
$ws.Range("G24").Value = "13"
$ws.Range("E25").Value = "5.45%"
$ws.Range("G25").Value = "13"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("G26").Value = "13"
$ws.Range("D27").Value = "0.0001690"
$ws.Range("E27").Value = "-12.78%"
$ws.Range("G27").Value = "13"
$ws.Range("G28").Value = "13"
$ws.Range("G29").Value = "13"
$ws.Range("G30").Value = "13"
$ws.Range("G31").Value = "13"
$ws.Range("G32").Value = "13"
$ws.Range("G33").Value = "13"
$ws.Range("G34").Value = "13"
$ws.Range("G35").Value = "13"
$ws.Range("G36").Value = "13"
$ws.Range("G37").Value = "13"
$ws.Range("G38").Value = "13"
$ws.Range("G39").Value = "13"
$ws.Range("D40").Value = "0.04089"
$ws.Range("E40").Value = "3.56%"
$ws.Range("G40").Value = "13"
$ws.Range("D41").Value = "0.006733"
$ws.Range("E41").Value = "61.60%"
$ws.Range("G41").Value = "13"
$ws.Range("D42").Value = "0.1164"
$ws.Range("E42").Value = "4.75%"
$ws.Range("G42").Value = "13"
$ws.Range("E43").Value = "4.72%"
$ws.Range("G43").Value = "13"
$ws.Range("D44").Value = "0.01247"
$ws.Range("E44").Value = "-10.33%"
$ws.Range("G44").Value = "13"
$ws.Range("D45").Value = "0.00005230"
$ws.Range("E45").Value = "2.31%"
$ws.Range("G45").Value = "13"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("G46").Value = "13"
$ws.Range("D47").Value = "2.377"
$ws.Range("E47").Value = "819.50%"
$ws.Range("G47").Value = "13"
$ws.Range("E48").Value = "6.15%"
$ws.Range("G48").Value = "13"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("G49").Value = "13"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("G50").Value = "13"
$ws.Range("G51").Value = "13"

# Restore the default (Normal) style so no stray number-format attributes
# remain on the cells after writing the text values above.
$dataRange.Style = "Normal"
